# Update "想去人数" (want-to-go count) figures in column F on the two
# sheets that carry the full listing data: "展览" (sheet 1) and
# "全部类型" (sheet 4). The "演出" and "本地生活" sheets don't contain
# these rows so they are left untouched.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (index 1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value  = 1158
$ws1.Range("F6").Value  = 14315
$ws1.Range("F7").Value  = 16424
$ws1.Range("F10").Value = 8
$ws1.Range("F12").Value = 199
$ws1.Range("F21").Value = 1253
$ws1.Range("F26").Value = 6660
$ws1.Range("F27").Value = 971
$ws1.Range("F28").Value = 4
$ws1.Range("F30").Value = 1116
$ws1.Range("F32").Value = 5730
$ws1.Range("F33").Value = 100
$ws1.Range("F34").Value = 144
$ws1.Range("F35").Value = 185
$ws1.Range("F36").Value = 4792

# --- Sheet "全部类型" (index 4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value  = 1158
$ws4.Range("F6").Value  = 14315
$ws4.Range("F7").Value  = 16424
$ws4.Range("F10").Value = 8
$ws4.Range("F12").Value = 199
$ws4.Range("F21").Value = 1253
$ws4.Range("F27").Value = 6660
$ws4.Range("F28").Value = 971
$ws4.Range("F29").Value = 4
$ws4.Range("F31").Value = 1116
$ws4.Range("F35").Value = 5730
$ws4.Range("F36").Value = 100
$ws4.Range("F37").Value = 144
$ws4.Range("F38").Value = 185
$ws4.Range("F39").Value = 4792
